# Updated cryptos list on Thu Feb  1 03:53:58 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.071.07"
$ws.Range("E2").Value = "  -2.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.259.33"
$ws.Range("E3").Value = "  -3.54%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "298.68"
$ws.Range("E5").Value = "  -2.68%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'94.10"
$ws.Range("E6").Value = "  -6.78%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.498"
$ws.Range("E7").Value = "  -2.44%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.80%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "32.99"
$ws.Range("E10").Value = "  -5.63%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -1.79%  "

# Row 12 - OKB
$ws.Range("D12").Value = "47.48"
$ws.Range("E12").Value = "  -8.81%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.41%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.91%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.610.87"
$ws.Range("E15").Value = "  -3.60%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "15.25"
$ws.Range("E16").Value = "  -4.13%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.254.77"
$ws.Range("E17").Value = "  -1.10%  "

# Row 18 - Polygon
$ws.Range("D18").Value = "0.774"
$ws.Range("E18").Value = "  -4.58%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "42.074.98"
$ws.Range("E19").Value = "  -1.85%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0892"

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -4.04%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "11.37"
$ws.Range("E22").Value = "  -2.68%  "

# Row 24 - BitcoinCash
$ws.Range("E24").Value = "  -1.52%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -5.17%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -4.28%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "23.72"
$ws.Range("E28").Value = "  -6.96%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'2.30"
$ws.Range("E29").Value = "  -1.14%  "

# Row 30 - Monero
$ws.Range("D30").Value = "167.15"
$ws.Range("E30").Value = "  +4.28%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").Value = "33.53"
$ws.Range("E31").Value = "  -4.55%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  -3.51%  "

# Row 33 - FirstDigitalUSD
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.07%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -3.87%  "

# Row 35 - WEMIXToken
$ws.Range("E35").Value = "  -5.75%  "

# Row 36 - was Hedera, now RenderToken
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.39"
$ws.Range("E36").Value = "  -6.06%  "

# Row 37 - was RenderToken, now Hedera
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0694"
$ws.Range("E37").Value = "  -4.74%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -6.40%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "15.93"
$ws.Range("E39").Value = "  -8.80%  "

# Row 40 - Kaspa
$ws.Range("D40").Value = "0.0991"
$ws.Range("E40").Value = "  -3.65%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  -3.61%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -8.73%  "

# Row 43 - ApeXProtocol
$ws.Range("E43").Value = "  -0.69%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.948.81"
$ws.Range("E44").Value = "  -3.67%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  -2.67%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "17.38"
$ws.Range("E46").Value = "  -7.28%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  -7.29%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  -5.18%  "

# Row 49 - HuobiToken
$ws.Range("E49").Value = "  -3.21%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.484.98"
$ws.Range("E50").Value = "  -2.94%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "52.16"
$ws.Range("E51").Value = "  -7.45%  "
